$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "Implemented turret death sound"
#
# Insert a new top-level bullet "TurretDeath" (with its two source
# links) right after the existing "HealthPickup" entry and its link,
# i.e. immediately before the "LowHealth" entry.
# ------------------------------------------------------------------

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

# Locate the paragraph that holds the "HealthPickup" hyperlink (the
# sub-bullet right under the "HealthPickup" heading) - the new
# material is inserted right after it.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq "HealthPickup") {
        $anchorIndex = $i + 1
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'HealthPickup' paragraph"
}

$pAnchor = $d.Paragraphs.Item($anchorIndex)

# 1) New top-level bullet: "TurretDeath"
$pAnchor.Range.InsertParagraphAfter() | Out-Null
$pTurretDeath = $d.Paragraphs.Item($anchorIndex + 1)
$pTurretDeath.Range.ListFormat.ListLevelNumber = 1
$pTurretDeath.Range.Text = "TurretDeath"

# 2) Sub-bullet hyperlink: steampunk machine sound
$pTurretDeath.Range.InsertParagraphAfter() | Out-Null
$pLink1 = $d.Paragraphs.Item($anchorIndex + 2)
$pLink1.Range.ListFormat.ListLevelNumber = 2
$url1 = "https://www.zapsplat.com/music/steampunk-machine-operating-clicks-gears-turning/"
$pLink1.Range.Text = $url1
$pLink1 = $d.Paragraphs.Item($anchorIndex + 2)
$link1Range = $d.Range($pLink1.Range.Start, $pLink1.Range.End - 1)
$d.Hyperlinks.Add($link1Range, $url1) | Out-Null

# 3) Sub-bullet hyperlink: grenade explosion sound
$pLink1 = $d.Paragraphs.Item($anchorIndex + 2)
$pLink1.Range.InsertParagraphAfter() | Out-Null
$pLink2 = $d.Paragraphs.Item($anchorIndex + 3)
$pLink2.Range.ListFormat.ListLevelNumber = 2
$url2 = "https://www.zapsplat.com/music/grenade-explosion-debris-blast/"
$pLink2.Range.Text = $url2
$pLink2 = $d.Paragraphs.Item($anchorIndex + 3)
$link2Range = $d.Range($pLink2.Range.Start, $pLink2.Range.End - 1)
$d.Hyperlinks.Add($link2Range, $url2) | Out-Null

Write-Host "Inserted TurretDeath entry with $($d.Hyperlinks.Count) total hyperlinks now in the document."
